# Update in ResetPassword and FYP_Test_Plan
# Adds two new test-case rows (TC-011 "Initialize App" and TC-012 "Reset
# Password") to the bottom of the "Test Plan" worksheet, below the existing
# Sign-up test-plan table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlCenter = -4108, xlLeft = -4131

# ---------------------------------------------------------------------
# Row 26 - TC-011 / Initialize App
# ---------------------------------------------------------------------
$ws.Range("A26").Value = 11
$ws.Range("B26").Value = "TC-011"
$ws.Range("C26").Value = "Initialize App"

$ws.Range("D26").Value = "S3 TOKEN GET NULL when start the app"
$ws.Range("D26").HorizontalAlignment = -4108
$ws.Range("D26").VerticalAlignment = -4108
$ws.Range("D26").WrapText = $true

$ws.Range("E26").Value = "1. Delete the previous account data in DynamoDB `n2. Delete the previous account data in Cognito `n3. Uninstall the existed App in phone `n4. Reinstall the app"
$ws.Range("E26").HorizontalAlignment = -4131
$ws.Range("E26").VerticalAlignment = -4108
$ws.Range("E26").WrapText = $true

$ws.Rows.Item(26).RowHeight = 100.8

# ---------------------------------------------------------------------
# Row 27 - TC-012 / Reset Password
# ---------------------------------------------------------------------
$ws.Range("A27").Value = 12
$ws.Range("B27").Value = "TC-012"
$ws.Range("C27").Value = "Reset Password"

$ws.Range("D27").Value = "Check whether the encrypted file by old password open correctly after change password"
$ws.Range("D27").HorizontalAlignment = -4108
$ws.Range("D27").VerticalAlignment = -4108
$ws.Range("D27").WrapText = $true

$ws.Range("E27").Value = "1. Create Account`n2. Create new area and new file`n3. Change Password via Reset Password`n4. Relogin the account with new password`n5. Open the created area and file and check whetehr it opens correctly"
$ws.Range("E27").HorizontalAlignment = -4131
$ws.Range("E27").VerticalAlignment = -4108
$ws.Range("E27").WrapText = $true

$ws.Rows.Item(27).RowHeight = 144

# ---------------------------------------------------------------------
# Update the view: scroll down to the newly-added rows and move the
# active selection to the next empty description cell.
# ---------------------------------------------------------------------
[void]$ws.Range("E28").Select()
